# scf_trip_time.xlsx — "Add files via upload" edit
#
# 1. Shared string "路程" -> "route_time" (used as the header in E1 of Sheet2)
# 2. Selection on Sheet2 moves from K15 to H27
# 3. Rows 30-44 on Sheet2 get a route_time (column E) value of 24;
#    row 44 additionally gets a trip_time (column D) value of 74 that it
#    was previously missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# 1. Rename the header text that used to read "路程" to "route_time".
$ws.Range("E1").Value = "route_time"

# 2. Add the missing trip_time value for row 44, then fill in route_time
#    (column E) with 24 for every data row from 30 through 44.
$ws.Range("D44").Value = 74

for ($r = 30; $r -le 44; $r++) {
    $ws.Cells.Item($r, 5).Value = 24
}

# 3. Move the active selection to H27 (was K15), keeping the same sheet active.
$ws.Activate() | Out-Null
$ws.Range("H27").Select() | Out-Null
